$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date in A1 (2024-04-24 -> 2024-05-24, serial 45406 -> 45436)
$ws.Range("A1").Value = 45436

# Update price cells
$ws.Range("D14").Value = 271.126
$ws.Range("D15").Value = 415.87
$ws.Range("D38").Value = 499.042
$ws.Range("D39").Value = 535.769

# Re-order the merged cell ranges by unmerging and re-merging in the
# desired order: A12:E12, A10:E10, A11:E11, A1:E1, A36:E36
$ws.Range("A12:E12").UnMerge()
$ws.Range("A10:E10").UnMerge()
$ws.Range("A11:E11").UnMerge()
$ws.Range("A1:E1").UnMerge()
$ws.Range("A36:E36").UnMerge()

$ws.Range("A12:E12").Merge()
$ws.Range("A10:E10").Merge()
$ws.Range("A11:E11").Merge()
$ws.Range("A1:E1").Merge()
$ws.Range("A36:E36").Merge()
